$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values per row, repulled data
$updates = @{
    2  = 2
    3  = 3
    5  = 3
    6  = -4
    7  = -2
    9  = 5
    10 = -4
    12 = -3
    13 = -5
    14 = 1
    15 = -5
    16 = -1
    17 = 1
    18 = 1
    19 = -5
    20 = 3
    21 = -3
    22 = 5
    24 = 1
    25 = -4
    26 = -2
    28 = 3
    29 = -3
    30 = 2
    31 = 3
    32 = 6
    34 = -2
    35 = -3
    36 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
